$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "video" flag column (F) for meeting001 and meeting003 to indicate
# images are now being grabbed via av instead of ImageMagick.
$ws.Range("F2").Value = 1
$ws.Range("F4").Value = 1

# Update the active selection on the sheet to reflect the new cursor position.
$ws.Range("F5").Select()
